# Add two new "history parent entity" modality rows to the value table.
#  1. A new "langue_sigle" / "ZZ" / "xyz" row inserted right after the
#     existing langue_sigle rows (before the oui_non rows), pushing the
#     later rows down by one.
#  2. A new "vide" / "missing" / "manquant" row appended at the end of
#     the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Insert the new row in the middle of the table (row 40) ----------------
# Inserting a whole row shifts rows 40:42 down to 41:43 and copies the
# formatting (style) from the row above, matching the target file.
$ws.Rows("40:40").Insert()

$ws.Range("A40").Value = "langue_sigle"
$ws.Range("B40").Value = "ZZ"
$ws.Range("C40").Value = "xyz"

# --- Append the new row at the end of the table (row 44) -------------------
$ws.Range("A44").Value = "vide"
$ws.Range("B44").Value = "missing"
$ws.Range("C44").Value = "manquant"

# --- Grow the table / autofilter range to cover the two new rows -----------
$lo.Resize($ws.Range("A1:C44"))

# --- Restore the selection shown in the saved file --------------------------
$ws.Range("C39").Select()
